# Update "Fruta, Agrícola del Norte S.A. de Arica - Pera" weekly price records
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44280
$ws.Range("M2").Value = 350
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 24500
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1361
$ws.Range("D3").Value = 44280
$ws.Range("K3").Value = 'Winter Nelis'
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("S3").Value = 1361
$ws.Range("D4").Value = 44336
$ws.Range("K4").Value = 'Winter Nelis'
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("S4").Value = 1194
$ws.Range("D5").Value = 44323
$ws.Range("K5").Value = 'Packham''s Triumph'
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("S5").Value = 861
$ws.Range("D6").Value = 44355
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 972
$ws.Range("D7").Value = 44355
$ws.Range("K7").Value = 'Winter Nelis'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 17500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 972
$ws.Range("D8").Value = 44371
$ws.Range("K8").Value = 'Packham''s Triumph'
$ws.Range("L8").Value = 'Calibre 90'
$ws.Range("M8").Value = 140
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17429
$ws.Range("Q8").Value = '$/caja 18 kilos embalada'
$ws.Range("S8").Value = 968
$ws.Range("D9").Value = 44371
$ws.Range("L9").Value = 'Calibre 80'
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 17500
$ws.Range("Q9").Value = '$/caja 18 kilos embalada'
$ws.Range("S9").Value = 972
$ws.Range("D10").Value = 44292
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 22500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 1250
$ws.Range("D11").Value = 44292
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 23000
$ws.Range("P11").Value = 22500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 1250
$ws.Range("D12").Value = 44421
$ws.Range("M12").Value = 270
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16500
$ws.Range("Q12").Value = '$/bandeja 18 kilos granel'
$ws.Range("S12").Value = 917
$ws.Range("D13").Value = 44421
$ws.Range("L13").Value = 'Segunda'
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 16500
$ws.Range("S13").Value = 917
$ws.Range("D14").Value = 44398
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 20500
$ws.Range("Q14").Value = '$/caja 20 kilos empedrada'
$ws.Range("S14").Value = 1025
$ws.Range("T14").Value = 20
$ws.Range("D15").Value = 44398
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20500
$ws.Range("Q15").Value = '$/caja 20 kilos empedrada'
$ws.Range("S15").Value = 1025
$ws.Range("T15").Value = 20
$ws.Range("D16").Value = 44341
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 17500
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 972
$ws.Range("T16").Value = 18
$ws.Range("D17").Value = 44313
$ws.Range("L17").Value = 'Tercera'
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15500
$ws.Range("Q17").Value = '$/bandeja 18 kilos granel'
$ws.Range("S17").Value = 861
$ws.Range("T17").Value = 18
$ws.Range("D18").Value = 44329
$ws.Range("M18").Value = 340
$ws.Range("N18").Value = 21000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 21500
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1194
$ws.Range("D19").Value = 44314
$ws.Range("K19").Value = 'Packham''s Triumph'
$ws.Range("R19").Value = 'Región de O''Higgins'

Write-Host "Applied weekly fruit/vegetable price updates"
